# Swap the colour palettes that live on the Slide Master's theme
# (ppt/theme/theme1.xml) and the Notes Master's theme (ppt/theme/theme2.xml).
# In the authored change the Slide Master ends up carrying the stock
# "Office" palette while the Notes Master ends up carrying the old
# "Integral" / "Red Violet" palette that the Slide Master used to have.

$p = $ppt.ActivePresentation

$masterScheme = $p.SlideMaster.Theme.ThemeColorScheme
$notesScheme  = $p.NotesMaster.Theme.ThemeColorScheme

$count = $masterScheme.Count

# Snapshot both palettes first so the swap doesn't clobber values we
# still need to read.
$masterRGB = @()
$notesRGB  = @()
for ($i = 1; $i -le $count; $i++) {
    $masterRGB += $masterScheme.Item($i).RGB
    $notesRGB  += $notesScheme.Item($i).RGB
}

for ($i = 1; $i -le $count; $i++) {
    $masterScheme.Item($i).RGB = $notesRGB[$i - 1]
    $notesScheme.Item($i).RGB  = $masterRGB[$i - 1]
}

# Point the three report tables (slides 14-16) at the other built-in
# table style.
$tableSlides = @(14, 15, 16)
foreach ($idx in $tableSlides) {
    $slide = $p.Slides.Item($idx)
    foreach ($shp in $slide.Shapes) {
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle("{26090D55-F4C2-48EF-B29D-A0A75C51BBED}")
        }
    }
}
